$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes existing rows 3..21 down to 4..22)
$ws.Rows("3:3").Insert()

# Populate the newly inserted row with the new IPO listing
$ws.Range("A3").Value = "아이스크림미디어(구.시공미디어)"
$ws.Range("B3").Value = "2024.07.31~08.06"
$ws.Range("C3").Value = "32,000~40,200"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 78720
$ws.Range("F3").Value = "삼성증권"

# The "SK증권스팩13호" row (now row 17 after the shift) gets its 확정공모가 (D) updated.
# Force text storage (matches the sibling cells in this column, e.g. D19/D20/D21) instead
# of letting the numeric-looking string be auto-converted to a number, then strip the
# quote-prefix style marker so the cell keeps the sheet's default (unstyled) format.
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2000"
$ws.Range("D17").ClearFormats()

# The last existing row (에이치브이엠, now pushed to row 22) drops off the bottom of the
# dataset range, so remove it to keep the sheet at rows 1-21
$ws.Rows("22:22").Delete()
